$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arduino-UltraSonic")
$ws.Activate()
